$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "30.265.32"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.862.04"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "236.38"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4709"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.98%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2896"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.20%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06542"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.19%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.98"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.07%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07952"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.99%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "97.62"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "1.872.70"
$ws.Range("E13").Value = "  +0.30%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.137"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.46%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6805"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "263.83"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -6.16%  "
$ws.Range("D17").Value = "30.267.51"
$ws.Range("E17").Value = "  +0.09%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.69"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +7.95%  "
$ws.Range("E19").Value = "  +0.07%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007476"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").Value = "2.116.16"
$ws.Range("E21").Value = "  +0.10%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.17%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.267"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -4.63%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.166"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "167.62"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.69%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.174"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.87"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.68%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.946"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.397"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.64%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.09902"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.18%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.344"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("E32").Value = "  -0.25%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.018"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.14%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04705"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.128"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7007"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.78%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.711"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.69%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01876"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.623"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.36%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.306"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.76%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "73.80"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.937"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4162"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "103.29"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "7.157"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.51%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "946.05"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.202"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.05%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "34.14"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "
